# Apply the "TEST_2" split-off edit to the HSK vocabulary workbook.
#
# What happened in the source edit:
#   1. The existing "TEST" sheet (char/pinyin/definition quiz data) was
#      duplicated to a new sheet named "TEST_2", placed right after "TEST".
#   2. A new row was appended to "TEST_2" (酋批 / nowei / no way), styled
#      like the sheet's existing plain (non-highlighted) rows.
#   3. The original "TEST" sheet was cleared out and repurposed as a small,
#      fresh two-column (word/definition) scratch/quiz sheet.
#
# Styles are copied (format-only) from the matching cells that already carry
# the right look, instead of being rebuilt property-by-property, so no
# stray/duplicate style entries are left behind.
#
# Shared-string insertion order matters for matching the target exactly:
# the new TEST content is typed first, then TEST_2's extra row - so we do
# the same here.

$wb = $excel.ActiveWorkbook
$testSheet = $wb.Worksheets.Item("TEST")

# ---------------------------------------------------------------------
# Step 1: duplicate "TEST" -> "TEST_2" (right after "TEST") BEFORE we
# touch/clear the original, so TEST_2 keeps all of TEST's current data
# and per-cell formatting intact.
# ---------------------------------------------------------------------
$testSheet.Copy([System.Reflection.Missing]::Value, $testSheet)
$test2 = $wb.Worksheets.Item("TEST (2)")
$test2.Name = "TEST_2"

$xlPasteFormats = -4122

# Grab references to a few still-intact TEST_2 cells whose formatting we
# want to reuse (header style, highlighted-row style, plain-row style).
$headerStyleSrc = $test2.Cells.Item(1,1)   # bold header look
$colAHighlightSrc = $test2.Cells.Item(2,1) # orange "char" column look
$colBHighlightSrc = $test2.Cells.Item(2,2) # orange "pinyin/definition" look
$colAPlainSrc = $test2.Cells.Item(8,1)     # plain "char" column look
$colBPlainSrc = $test2.Cells.Item(8,2)     # plain "pinyin/definition" look

# ---------------------------------------------------------------------
# Step 2: clear the original "TEST" sheet and give it fresh content -
# a little word/definition scratch list.
# ---------------------------------------------------------------------
$testSheet.Cells.Clear()

$headerStyleSrc.Copy()
$testSheet.Cells.Item(1,1).PasteSpecial($xlPasteFormats)
$testSheet.Cells.Item(1,1).Value = "word"

$headerStyleSrc.Copy()
$testSheet.Cells.Item(1,2).PasteSpecial($xlPasteFormats)
$testSheet.Cells.Item(1,2).Value = "definition"

$newTestData = @(
  @("baba", "baba is…"),
  @("potato", "potato means…"),
  @("a", "a is just a letter"),
  @("f", "f is also just a letter")
)

$r = 2
foreach ($pair in $newTestData) {
  $colAHighlightSrc.Copy()
  $testSheet.Cells.Item($r,1).PasteSpecial($xlPasteFormats)
  $testSheet.Cells.Item($r,1).Value = $pair[0]

  $colBHighlightSrc.Copy()
  $testSheet.Cells.Item($r,2).PasteSpecial($xlPasteFormats)
  $testSheet.Cells.Item($r,2).Value = $pair[1]

  $r++
}

for ($rr = 1; $rr -le 15; $rr++) {
  $testSheet.Rows.Item($rr).RowHeight = 20.9
}

$testSheet.Columns.Item(1).ColumnWidth = 8.1
$testSheet.Columns.Item(2).ColumnWidth = 25.85

# ---------------------------------------------------------------------
# Step 3: append the extra vocabulary row to "TEST_2" (plain style,
# matching the sheet's existing non-highlighted rows, e.g. row 8).
# ---------------------------------------------------------------------
$colAPlainSrc.Copy()
$test2.Cells.Item(9,1).PasteSpecial($xlPasteFormats)
$test2.Cells.Item(9,1).Value = "酋批"

$colBPlainSrc.Copy()
$test2.Cells.Item(9,2).PasteSpecial($xlPasteFormats)
$test2.Cells.Item(9,2).Value = "nowei"

$colBPlainSrc.Copy()
$test2.Cells.Item(9,3).PasteSpecial($xlPasteFormats)
$test2.Cells.Item(9,3).Value = "no way"

for ($rr = 9; $rr -le 15; $rr++) {
  $test2.Rows.Item($rr).RowHeight = 20.9
}

$test2.Range("B13").Select()

# ---------------------------------------------------------------------
# Step 4: leave "TEST" as the active/selected tab (matches the
# workbook's original activeTab, which keeps pointing at "TEST").
# ---------------------------------------------------------------------
$testSheet.Activate()
$testSheet.Range("E15").Select()
